$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with latest values.
# Some Price values are plain decimal numbers (e.g. "576.45"); force them to remain
# text cells (matching the original inlineStr formatting) by temporarily applying a
# text number format, then restoring the default "Normal" style so no visible
# formatting change is introduced.
$ws.Range('D2').Value = '61.896.23'
$ws.Range('E2').Value = '  +3.00%  '
$ws.Range('D3').Value = '3.401.72'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.06%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.46'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('E10').Value = '  +6.96%  '
$ws.Range('E11').Value = '  +4.21%  '
$ws.Range('D12').Value = '3.984.14'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000177'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.98%  '
$ws.Range('D15').Value = '3.401.23'
$ws.Range('E15').Value = '  +1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.47%  '
$ws.Range('D17').Value = '61.953.41'
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('E18').Value = '  +5.75%  '
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('E20').Value = '  +4.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '388.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.58%  '
$ws.Range('E22').Value = '  +1.97%  '
$ws.Range('D23').Value = '3.545.78'
$ws.Range('E23').Value = '  +2.11%  '
$ws.Range('E24').Value = '  +14.75%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.08%  '
$ws.Range('E28').Value = '  -4.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  +4.40%  '
$ws.Range('E31').Value = '  +4.83%  '
$ws.Range('E32').Value = '  +2.37%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').Value = '3.433.46'
$ws.Range('E34').Value = '  +1.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.52'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.41'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.98'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.55'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '163.12'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.49%  '
$ws.Range('E40').Value = '  +2.46%  '
$ws.Range('E41').Value = '  +13.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.785'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.74%  '
$ws.Range('E43').Value = '  +4.42%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.86'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.46%  '
$ws.Range('E48').Value = '  +2.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.30'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.19%  '
$ws.Range('D50').Value = '2.373.67'
$ws.Range('E50').Value = '  +8.71%  '
$ws.Range('E51').Value = '  +4.89%  '
